$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A6 used to hold the numeric AHV-Nr 7560000000004. Replace it (and the row)
# with string-formatted values, and add two more rows that exercise other
# ahv_format edge cases the pseudonymizer needs to detect.

$ws.Range("A6").Value = "756.0000.000.004"
$ws.Range("B6").Value = "Dachs"
$ws.Range("C6").Value = "Dario"

$ws.Range("A7").Value = "756AB00000004"
$ws.Range("B7").Value = "Eris"
$ws.Range("C7").Value = "Elsa"

$ws.Range("A8").Value = "7230000XYZ"
$ws.Range("B8").Value = "Fichter"
$ws.Range("C8").Value = "Fiona"

# Keep the integer-style number format on A7 (same as A6/A2:A4), but A8
# should fall back to the default (unstyled) cell format.
$ws.Range("A7").NumberFormat = "0"

# Move the active selection to A6 to match where the edits were made.
$ws.Range("A6").Select()
